# ============================================================================
# Edit script: adds "Player Info" and "ODI Batting Extra" sheets, and updates
# the existing "ODI Batting" sheet (renames MATCH_CARD_LINK -> MATCH_CODE,
# converts the scorecard URLs to bare match codes, and clears a handful of
# blank INNING_NUMBER cells).
# ============================================================================

$wb = $excel.ActiveWorkbook
$battingSheet = $wb.Worksheets.Item(1)

# ----------------------------------------------------------------------------
# 1. Update the existing "ODI Batting" sheet
# ----------------------------------------------------------------------------

# Rename the MATCH_CARD_LINK header to MATCH_CODE
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

# Convert every MATCH_CARD_LINK URL in column D into just the bare match code
$battingSheet.Range("D2:D84").NumberFormat = "@"
for ($r = 2; $r -le 84; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $url = $cell.Value()
    $idx = $url.LastIndexOf("=")
    $code = $url.Substring($idx + 1)
    $cell.Value = $code
}

# Clear the handful of blank INNING_NUMBER (column B) cells belonging to
# "did not bat" rows so that no empty cell element remains
$blankInningRows = @(12, 20, 22, 23, 71, 80, 81)
foreach ($r in $blankInningRows) {
    $battingSheet.Cells.Item($r, 2).ClearContents()
}

# ----------------------------------------------------------------------------
# 2. Add the "Player Info" sheet (before "ODI Batting")
# ----------------------------------------------------------------------------

$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1:D2").NumberFormat = "@"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 0; $c -lt $piHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c + 1)
    $cell.Value = $piHeaders[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$piRow = @("4221", "Matthew Henry Cross", "Right Handed", "Does Not Bowl | Unknown")
for ($c = 0; $c -lt $piRow.Length; $c++) {
    $playerInfo.Cells.Item(2, $c + 1).Value = $piRow[$c]
}

# ----------------------------------------------------------------------------
# 3. Add the "ODI Batting Extra" sheet (after "ODI Batting")
# ----------------------------------------------------------------------------

$battingSheet2 = $wb.Worksheets.Item("ODI Batting")
$battingExtra = $wb.Worksheets.Add($null, $battingSheet2)
$battingExtra.Name = "ODI Batting Extra"

$battingExtra.Range("A1:A21").NumberFormat = "@"
$battingExtra.Range("B1").NumberFormat = "@"
$battingExtra.Range("C1:F21").NumberFormat = "@"

$beHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 0; $c -lt $beHeaders.Length; $c++) {
    $cell = $battingExtra.Cells.Item(1, $c + 1)
    $cell.Value = $beHeaders[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$beRows = @(
    @("4576", 1, "3", "2", "24.67%", "NO"),
    @("4578", 1, "0", "0", "0.92%", "NO"),
    @("4581", 1, "1", "0", "2.92%", "NO"),
    @("4604", 2, "3", "1", "10.85%", "NO"),
    @("4610", 2, "0", "0", "1.39%", "NO"),
    @("4612", 2, "4", "0", "13.24%", "NO"),
    @("4617", 6, $null, $null, $null, "NO"),
    @("4625", 5, "7", "0", "17.32%", "NO"),
    @("4629", 5, "1", "0", "9.54%", "NO"),
    @("4631", $null, $null, $null, $null, "NO"),
    @("4632", 5, "4", "1", "33.46%", "NO"),
    @("4635", $null, $null, $null, $null, "NO"),
    @("4677", 5, "3", "0", "8.92%", "NO"),
    @("4681", 5, "1", "0", "9.09%", "NO"),
    @("4680", 5, "1", "0", "8.17%", "NO"),
    @("4684", 5, $null, $null, $null, "NO"),
    @("4702", 4, $null, $null, $null, "NO"),
    @("4703", 4, "0", "0", "0.73%", "NO"),
    @("4705", 5, "1", "1", "7.24%", "NO"),
    @("4706", $null, $null, $null, $null, "NO")
)

for ($i = 0; $i -lt $beRows.Length; $i++) {
    $rowVals = $beRows[$i]
    $r = $i + 2
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $val = $rowVals[$c]
        if ($val -ne $null) {
            $battingExtra.Cells.Item($r, $c + 1).Value = $val
        }
    }
}

Write-Host "Edit complete"
